# Update HotStock_Top20 data (columns A:C, rows 2-21) to reflect the
# latest "hot stock" rankings from 财联社 (A), 东方财富 (B) and 同花顺 (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("中际旭创", "寒武纪-U", "卧龙电驱"),
    @("卧龙电驱", "岩山科技", "农业银行"),
    @("岩山科技", "卧龙电驱", "三维通信"),
    @("寒武纪-U", "长城军工", "华胜天成"),
    @("春兴精工", "太平洋", "岩山科技"),
    @("阳光电源", "吉视传媒", "寒武纪"),
    @("三维通信", "三维通信", "太平洋"),
    @("太平洋", "中际旭创", "工业富联"),
    @("景兴纸业", "春兴精工", "北方稀土"),
    @("国轩高科", "景兴纸业", "供销大集"),
    @("吉视传媒", "供销大集", "新易盛"),
    @("利欧股份", "利欧股份", "利欧股份"),
    @("新易盛",   "新易盛",   "中际旭创"),
    @("上能电气", "步步高",   "国轩高科"),
    @("供销大集", "上能电气", "金财互联"),
    @("亿纬锂能", "工业富联", "山子高科"),
    @("长城军工", "东方财富", "领益智造"),
    @("工业富联", "农业银行", "东方财富"),
    @("华胜天成", "阳光电源", "金龙羽"),
    @("北方稀土", "北方稀土", "春兴精工")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
